# amazonSteps.xlsx update:
# - replace the search/add-to-cart flow (old rows 10-18) with a streamlined
#   proceed-to-checkout / deliver-to-address flow (new rows 10-14)
# - refresh a handful of earlier locator strings (rows 3-9)
# - move the cart hyperlink from the old B14 "Proceed to checkout" cell to
#   the new B10 "goto cart" cell
# - update the active-cell selection to B11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. drop all existing hyperlinks up front so row surgery below can't ---
# --- leave stale/duplicated hyperlink entries behind                    ---
$ws.Range("A1").Hyperlinks.Delete()

# --- 2. remove the trailing rows (old 15-18) that no longer exist -----------
$ws.Rows("15:18").Delete()

# --- 3. copy B14's current "Hyperlink" cell style onto B10 before we --------
# --- overwrite B14's own content/value (both end up styled as a link) -------
$ws.Range("B14").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B10").HorizontalAlignment = -4131
$ws.Range("B10").VerticalAlignment = -4160

# --- 4. clear out the old values for rows 3-14 (contents only, keep style) --
$ws.Range("A3:E14").ClearContents()

# --- 5. write the new step data ---------------------------------------------
$ws.Range("A2").Value = "goto"
$ws.Range("B2").Value = "https://www.amazon.com/"
$ws.Range("D2").Value = 1000
$ws.Range("E2").Value = 3000

$ws.Range("A3").Value = "waitfortext"
$ws.Range("B3").Value = "Hello, Sign in"

$ws.Range("A4").Value = "click"
$ws.Range("B4").Value = "Hello, Sign in"
$ws.Range("D4").Value = 1000
$ws.Range("E4").Value = 3000

$ws.Range("A5").Value = "waitfortext"
$ws.Range("B5").Value = "Continue"

$ws.Range("A6").Value = "type"
$ws.Range("B6").Value = "email input field"
$ws.Range("C6").Value = "weavernormar@gmail.com"
$ws.Range("D6").Value = 1000
$ws.Range("E6").Value = 2000

$ws.Range("A7").Value = "click"
$ws.Range("B7").Value = "Continue button on page"
$ws.Range("D7").Value = 1000
$ws.Range("E7").Value = 3000

$ws.Range("A8").Value = "type"
$ws.Range("B8").Value = "password input field in span"
$ws.Range("C8").Value = "Welcome@123456"
$ws.Range("D8").Value = 1000
$ws.Range("E8").Value = 1000

$ws.Range("A9").Value = "click"
$ws.Range("B9").Value = "signin button on page"
$ws.Range("D9").Value = 1000
$ws.Range("E9").Value = 3000

$ws.Range("A10").Value = "goto"
$ws.Range("B10").Value = "https://www.amazon.com/gp/cart/view.html?ref_=nav_cart"
$ws.Range("D10").Value = 1000
$ws.Range("E10").Value = 7000

$ws.Range("A11").Value = "clickloc"
$ws.Range("B11").Value = "input[name=""proceedToRetailCheckout""]"
$ws.Range("D11").Value = 1000
$ws.Range("E11").Value = 7000

$ws.Range("A12").Value = "pressto"
$ws.Range("B12").Value = "Deliver to this address"
$ws.Range("C12").Value = "first"
$ws.Range("D12").Value = 1000
$ws.Range("E12").Value = 7000

$ws.Range("A13").Value = "assert"
$ws.Range("B13").Value = "h2#deliver-to-customer-text"
$ws.Range("C13").Value = "Delivering to Normar Weaver"
$ws.Range("D13").Value = 1000

# row 14 stays present but empty (B14 keeps its inherited Hyperlink style)

# --- 6. re-create the three hyperlinks -------------------------------------
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:weavernormar@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C8"), "mailto:Welcome@123456")
$ws.Hyperlinks.Add($ws.Range("B10"), "https://www.amazon.com/gp/cart/view.html?ref_=nav_cart")

# --- 7. restore the selection to B11 ----------------------------------------
$ws.Range("B11").Select()
